$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add P1=14, Q1=15, matching the bold/bordered style used by the rest of row 1 ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2..25: swap I<->K and M<->O, and append P=2, Q=2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O (was 2)
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
